# Level 2 framework changes
# - Update B7 value on "WeatherAPITestParameters" sheet from 21.14 to 17.38
# - Move the active sheet/tab selection from "CityNameCountryParameters" (A4)
#   to "WeatherAPITestParameters", with the in-sheet selection moved from
#   A12:F17 to the single cell B7.

$wb = $excel.ActiveWorkbook

$wsParams = $wb.Worksheets.Item("WeatherAPITestParameters")
$wsCityName = $wb.Worksheets.Item("CityNameCountryParameters")

# Update the data value.
$wsParams.Cells.Item(7, 2).Value = 17.38

# Make sure the previously active sheet's selection is left at A4 (unchanged),
# then switch activation to the WeatherAPITestParameters sheet and select B7.
$wsCityName.Activate()
$wsCityName.Range("A4").Select()

$wsParams.Activate()
$wsParams.Range("B7").Select()
